$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(152).Insert()

$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 45086
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = 100112001
$ws.Range("G152").Value = "Berenjena"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 300
$ws.Range("K152").Value = 9000
$ws.Range("L152").Value = 10000
$ws.Range("M152").Value = 9500
$ws.Range("N152").Value = "$/caja 50 unidades"
$ws.Range("O152").Value = "Región de Arica y Parinacota"
$ws.Range("P152").Value = 190
$ws.Range("Q152").Value = 50
$ws.Range("R152").Value = "Hortaliza"
